$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,  "T.J. McConnell",       "PG",       "Indiana Pacers"),
    @(3,  "Jared McCain",         "PG,SG",    "Philadelphia 76ers"),
    @(4,  "Derrick White",        "PG,SG",    "Boston Celtics"),
    @(5,  "Cam Thomas",           "SG,SF",    "Brooklyn Nets"),
    @(6,  "Cameron Johnson",      "SF,PF",    "Brooklyn Nets"),
    @(7,  "Julius Randle",        "PF",       "Minnesota Timberwolves"),
    @(8,  "Bilal Coulibaly",      "SG,SF",    "Washington Wizards"),
    @(9,  "Bam Adebayo",          "C",        "Miami Heat"),
    @(10, "LaMelo Ball",          "PG,SG",    "Charlotte Hornets"),
    @(11, "Damian Lillard",       "PG",       "Milwaukee Bucks"),
    @(12, "Isaiah Hartenstein",   "C",        "Oklahoma City Thunder"),
    @(13, "Cody Martin",          "SG,SF",    "Charlotte Hornets"),
    @(14, "Brandon Miller",       "SG,SF",    "Charlotte Hornets"),
    @(15, "Anthony Davis",        "PF,C",     "Los Angeles Lakers"),
    @(16, "Brandon Ingram",       "SG,SF,PF", "New Orleans Pelicans")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
